$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "59.60") are not coerced to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '37.323.31'
$ws.Range('E2').Value = '  +4.18%  '

$ws.Range('D3').Value = '2.041.84'
$ws.Range('E3').Value = '  +2.56%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '251.67'
$ws.Range('E5').Value = '  +2.39%  '

$ws.Range('D6').Value = '0.649'
$ws.Range('E6').Value = '  +1.52%  '

$ws.Range('D7').Value = '65.16'
$ws.Range('E7').Value = '  +9.35%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').Value = '0.401'
$ws.Range('E9').Value = '  +9.39%  '

$ws.Range('D10').Value = '59.60'
$ws.Range('E10').Value = '  +1.08%  '

$ws.Range('D11').Value = '0.0786'
$ws.Range('E11').Value = '  +6.17%  '

$ws.Range('D12').Value = '0.104'
$ws.Range('E12').Value = '  -0.46%  '

$ws.Range('D13').Value = '0.917'
$ws.Range('E13').Value = '  -3.22%  '

$ws.Range('D14').Value = '23.16'
$ws.Range('E14').Value = '  +19.99%  '

$ws.Range('D15').Value = '14.74'
$ws.Range('E15').Value = '  -0.03%  '

$ws.Range('D16').Value = '2.341.33'
$ws.Range('E16').Value = '  +2.67%  '

$ws.Range('D17').Value = '5.70'
$ws.Range('E17').Value = '  +6.90%  '

$ws.Range('D18').Value = '2.044.43'
$ws.Range('E18').Value = '  +2.71%  '

$ws.Range('D19').Value = '37.202.64'
$ws.Range('E19').Value = '  +4.08%  '

$ws.Range('D20').Value = '73.41'
$ws.Range('E20').Value = '  +2.13%  '

$ws.Range('D21').Value = '0.0₃0879'
$ws.Range('E21').Value = '  +3.43%  '

$ws.Range('D22').Value = '5.49'
$ws.Range('E22').Value = '  +5.30%  '

$ws.Range('D23').Value = '238.92'
$ws.Range('E23').Value = '  +2.27%  '

$ws.Range('E24').Value = '  -0.09%  '

$ws.Range('D25').Value = '2.60'
$ws.Range('E25').Value = '  +0.11%  '

$ws.Range('D26').Value = '2.37'
$ws.Range('E26').Value = '  +3.98%  '

$ws.Range('D27').Value = '10.00'
$ws.Range('E27').Value = '  +4.69%  '

$ws.Range('D28').Value = '161.27'
$ws.Range('E28').Value = '  -2.27%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '19.97'
$ws.Range('E29').Value = '  +3.17%  '

$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.133'
$ws.Range('E30').Value = '  +36.09%  '

$ws.Range('D31').Value = '0.122'
$ws.Range('E31').Value = '  +2.23%  '

$ws.Range('D32').Value = '5.14'
$ws.Range('E32').Value = '  +4.27%  '

$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  +4.63%  '

$ws.Range('D34').Value = '0.0626'
$ws.Range('E34').Value = '  +3.81%  '

$ws.Range('D35').Value = '4.67'
$ws.Range('E35').Value = '  +6.03%  '

$ws.Range('D36').Value = '2.37'
$ws.Range('E36').Value = '  -3.93%  '

$ws.Range('D37').Value = '6.30'
$ws.Range('E37').Value = '  +10.13%  '

$ws.Range('E38').Value = '  -0.05%  '

$ws.Range('E39').Value = '  +2.37%  '

$ws.Range('D40').Value = '2.97'
$ws.Range('E40').Value = '  +27.49%  '

$ws.Range('D41').Value = '0.101'
$ws.Range('E41').Value = '  +7.58%  '

$ws.Range('D42').Value = '1.27'
$ws.Range('E42').Value = '  +2.31%  '

$ws.Range('E43').Value = '  +6.76%  '

$ws.Range('D44').Value = '17.44'
$ws.Range('E44').Value = '  +6.17%  '

$ws.Range('D45').Value = '1.15'
$ws.Range('E45').Value = '  +4.60%  '

$ws.Range('D46').Value = '0.0218'
$ws.Range('E46').Value = '  +1.89%  '

$ws.Range('D47').Value = '95.72'
$ws.Range('E47').Value = '  +1.60%  '

$ws.Range('D48').Value = '7.83'
$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('D49').Value = '1.395.35'
$ws.Range('E49').Value = '  +1.77%  '

$ws.Range('D50').Value = '2.91'
$ws.Range('E50').Value = '  +0.50%  '

$ws.Range('D51').Value = '47.18'
$ws.Range('E51').Value = '  +0.56%  '

# Restore column D to the default (Normal) style so no stray number-format styling is introduced.
$dRange.Style = "Normal"